$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("APISuite")

# --- Row 2: uri / /booking for "Test Case 1 - Create Booking Postive scenario" ---
# Order matters for shared-string append order: uri, /booking, name, (reuse name), Request Body, json
$ws.Range("C2").Value = "uri"
$ws.Range("D2").Value = "/booking"
$ws.Range("B2").Value = "Test Case 1 - Create Booking Postive scenario"

# --- Row 3: Request Body + JSON payload ---
$ws.Range("B3").Value = "Test Case 1 - Create Booking Postive scenario"
$ws.Range("C3").Value = "Request Body"
$json = "{`n    ""bookingid"": 25,`n    ""booking"": {`n        ""firstname"": ""Jim"",`n        ""lastname"": ""Brown"",`n        ""totalprice"": 111,`n        ""depositpaid"": true,`n        ""bookingdates"": {`n            ""checkin"": ""2018-01-01"",`n            ""checkout"": ""2019-01-01""`n        },`n        ""additionalneeds"": ""Breakfast""`n    }`n}"
$ws.Range("D3").Value = $json

# Row 3 formatting: vertical-center the A/B/C cells, wrap + shrink font on the body cell
# (order controls the append order of newly-created cell styles)
$ws.Range("B3").VerticalAlignment = -4108
$ws.Range("C3").VerticalAlignment = -4108
$ws.Range("D3").WrapText = $true
$ws.Range("D3").Font.Size = 10
$ws.Range("A3").VerticalAlignment = -4108

$ws.Rows.Item(3).RowHeight = 180

# --- Page setup ---
$ws.PageSetup.Orientation = 1

# --- Selection / active sheet state ---
$ws1 = $wb.Worksheets.Item("UISuite")
$ws1.Range("B13").Select()

$ws.Activate()
$ws.Range("B3").Select()
